# Update the "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect newly generated output data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Row number -> new value for column F
$updates = @{
    5  = 40
    6  = 218
    10 = 5648
    11 = 5011
    15 = 57
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}

$wb.Save()
